$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A84").Value = "'2023-06-28"
$ws.Range("B84").Value = "'10:31:25"
$ws.Range("C84").Value = "'Wednesday"
$ws.Range("D84").Value = "'26"
$ws.Range("A84:D84").ClearFormats()

$ws.Range("E84").Value = 122902
$ws.Range("F84").Value = 134460
$ws.Range("G84").Value = 163476
$ws.Range("H84").Value = 134027
$ws.Range("I84").Value = 177199
$ws.Range("J84").Value = 114782
$ws.Range("K84").Value = 203875
$ws.Range("L84").Value = 226493
$ws.Range("M84").Value = 176128
$ws.Range("N84").Value = 104415
$ws.Range("O84").Value = 39715
$ws.Range("P84").Value = 33757
$ws.Range("Q84").Value = 52334
$ws.Range("R84").Value = -1
$ws.Range("S84").Value = 35606
$ws.Range("T84").Value = -1
